$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper function: set a cell's value as literal text, even when the text
# looks numeric (e.g. "0.08"), while preserving the destination cell's
# existing style/number-format. Assigning a formula that evaluates to a
# text string (="0.08") forces a text result; copying the cell and pasting
# values-only back onto itself then bakes that text result in as a plain
# (non-formula) value, without disturbing the cell's NumberFormat/style.
function Set-TextValue($destCell, [string]$text) {
    $escaped = $text -replace '"', '""'
    $destCell.Formula = '="' + $escaped + '"'
    $destCell.Copy()
    $destCell.PasteSpecial(-4163, 0, $false, $false)
}

# Row 3 and Row 4 swap their Material / Quantity / Inventoryitem values
# (Pantone-1 - UV / 0.14 / UV-PMS-UV-PMS  <->  Black - UV / 0.08 / UV-4 CP...)
$ws.Cells.Item(3, 4).Value = "Black - UV"
Set-TextValue $ws.Cells.Item(3, 5) "0.08"
$ws.Cells.Item(3, 8).Value = "UV-4 CP - UV- 4 Color Process"

$ws.Cells.Item(4, 4).Value = "Pantone-1 - UV"
Set-TextValue $ws.Cells.Item(4, 5) "0.14"
$ws.Cells.Item(4, 8).Value = "UV-PMS - UV-PMS"

# Row 6 and Row 7 swap their Element / Process / Material / Quantity /
# Inventoryitem values
$ws.Cells.Item(6, 2).Value = "Plate - Outer wrap  2p"
$ws.Cells.Item(6, 3).Value = "Plate burn"
$ws.Cells.Item(6, 4).Value = "CTP 40"""
Set-TextValue $ws.Cells.Item(6, 5) "2.00"
$ws.Cells.Item(6, 8).Value = "CTP40 - CTP Plate 40"""

$ws.Cells.Item(7, 2).Value = "Outer wrap  2p"
$ws.Cells.Item(7, 3).Value = "Print F (Varnish 1x0)"
$ws.Cells.Item(7, 4).Value = "Cyrel Plate"
Set-TextValue $ws.Cells.Item(7, 5) "1.00"
$ws.Cells.Item(7, 8).Value = "450 - ***40"" Cyrel Plate for Press_________***"

# Row 8 and Row 9 swap their Element / Process / Material / Quantity /
# Measurementunit / Inventoryitem values
$ws.Cells.Item(8, 2).Value = "Liner  2p"
$ws.Cells.Item(8, 3).Value = "Cut for Press"
$ws.Cells.Item(8, 4).Value = "Special Order Coated Book Non FSC 100 24 x 36"" 400 ppi"
Set-TextValue $ws.Cells.Item(8, 5) "150.00"
$ws.Cells.Item(8, 6).Value = "Sht."
$ws.Cells.Item(8, 8).Value = "429 - Offset Stock cost ea for Jobs Only-Use  for shipping too."

$ws.Cells.Item(9, 2).Value = "Outer wrap  2p"
$ws.Cells.Item(9, 3).Value = "Print F 2x0"
$ws.Cells.Item(9, 4).Value = "Sterling Ultra C1S Gloss Verso Non FSC 80# 19 x 25"" 500 ppi"
Set-TextValue $ws.Cells.Item(9, 5) "400.00"
$ws.Cells.Item(9, 6).Value = "Sht"
$ws.Cells.Item(9, 8).Value = "430 - Offset Stock Per 1000 cost for Jobs Only"
